$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 381, pushing the existing rows 381..463 down to 383..465
$ws.Rows.Item(381).Resize(2).Insert()

# Populate the two newly inserted rows (381 and 382) with the new Albahaca price records.

# Row 381
$ws.Cells.Item(381, 1).Value = 9
$ws.Cells.Item(381, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(381, 3).Value = "Metropolitana"
$ws.Cells.Item(381, 4).Value = 44932
$ws.Cells.Item(381, 5).Value = 13
$ws.Cells.Item(381, 6).Value = 100112052
$ws.Cells.Item(381, 7).Value = "Albahaca"
$ws.Cells.Item(381, 8).Value = "Sin especificar"
$ws.Cells.Item(381, 9).Value = "Primera"
$ws.Cells.Item(381, 10).Value = 340
$ws.Cells.Item(381, 11).Value = 4000
$ws.Cells.Item(381, 12).Value = 5000
$ws.Cells.Item(381, 13).Value = 4500
$ws.Cells.Item(381, 14).Value = "$/docena de matas"
$ws.Cells.Item(381, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(381, 16).Value = 750
$ws.Cells.Item(381, 17).Value = 6
$ws.Cells.Item(381, 18).Value = "Hortaliza"

# Row 382
$ws.Cells.Item(382, 1).Value = 9
$ws.Cells.Item(382, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(382, 3).Value = "Metropolitana"
$ws.Cells.Item(382, 4).Value = 44932
$ws.Cells.Item(382, 5).Value = 13
$ws.Cells.Item(382, 6).Value = 100112052
$ws.Cells.Item(382, 7).Value = "Albahaca"
$ws.Cells.Item(382, 8).Value = "Sin especificar"
$ws.Cells.Item(382, 9).Value = "Primera"
$ws.Cells.Item(382, 10).Value = 430
$ws.Cells.Item(382, 11).Value = 4000
$ws.Cells.Item(382, 12).Value = 5000
$ws.Cells.Item(382, 13).Value = 4500
$ws.Cells.Item(382, 14).Value = "$/docena de matas"
$ws.Cells.Item(382, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(382, 16).Value = 750
$ws.Cells.Item(382, 17).Value = 6
$ws.Cells.Item(382, 18).Value = "Hortaliza"
